{"js": "// Update benchmark stats table: refresh the summary column values and\n// collapse the last three multi-tab rows down to a single re-computed value.\nconst table = context.document.body.tables.getFirst();\n\n// Map of 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"1923\",\n  4: \"0.00002\",\n  5: \"0.00069\",\n  6: \"0.00017\",\n  7: \"0.00005\",\n  8: \"0.00028\",\n  9: \"0.00038\",\n  10: \"0.00048\",\n  11: \"0.38744\",\n  43: \"99.74\",\n  44: \"0.39\",\n  45: \"146\",\n};\n\nfor (const [rowIdx, newText] of Object.entries(updates)) {\n  const cell = table.getCell(Number(rowIdx), 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table: refresh the summary column values and\n# collapse the last three multi-tab rows down to a single re-computed value.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# 1-based row -> new cell text (table has a single column).\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"1923\"\n    5  = \"0.00002\"\n    6  = \"0.00069\"\n    7  = \"0.00017\"\n    8  = \"0.00005\"\n    9  = \"0.00028\"\n    10 = \"0.00038\"\n    11 = \"0.00048\"\n    12 = \"0.38744\"\n    44 = \"99.74\"\n    45 = \"0.39\"\n    46 = \"146\"\n}\n\nforeach ($row in $updates.Keys) {\n    $t.Cell($row, 1).Range.Text = $updates[$row]\n}\n"}
